$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 790.2632
$ws.Range("I19").Value = 216.5
$ws.Range("K19").Value = 216.5
$ws.Range("M19").Value = -41.5
$ws.Range("H30").Value = 980
$ws.Range("J30").Value = 980
$ws.Range("L30").Value = 2940
$ws.Range("N30").Value = -3142
$ws.Range("H88").Value = 1374708.1
$ws.Range("I88").Value = 800
$ws.Range("J88").Value = 1546446.6
$ws.Range("K88").Value = 800
$ws.Range("L88").Value = 1546446.6
$ws.Range("M88").Value = -394
$ws.Range("N88").Value = -1547258.6
$ws.Range("H91").Value = 1374708.1
$ws.Range("I91").Value = 800
$ws.Range("J91").Value = 1546446.6
$ws.Range("K91").Value = 800
$ws.Range("L91").Value = 1546446.6
$ws.Range("M91").Value = 604
$ws.Range("N91").Value = -1549254.6
$ws.Range("H98").Value = 3558.3333
$ws.Range("I98").Value = 3686.25
$ws.Range("K98").Value = 3686.25
$ws.Range("M98").Value = -2188.25
$ws.Range("H122").Value = 3558.3333
$ws.Range("I122").Value = 3686.25
$ws.Range("K122").Value = 11058.75
$ws.Range("M122").Value = -8608.75
$ws.Range("H138").Value = 1414.62
$ws.Range("I138").Value = 898.55884
$ws.Range("J138").Value = 1680.4697
$ws.Range("K138").Value = 2695.67652
$ws.Range("L138").Value = 5041.409100000001
$ws.Range("M138").Value = 2444.32348
$ws.Range("N138").Value = -15321.4091
$ws.Range("H141").Value = 2091.4443
$ws.Range("I141").Value = 965
$ws.Range("K141").Value = 2895
$ws.Range("M141").Value = 2285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4268.979
$ws.Range("I32").Value = 3904.5813
$ws.Range("K32").Value = 3904.5813
$ws.Range("M32").Value = -3617.5813
$ws.Range("H61").Value = 641.3125
$ws.Range("I61").Value = 641.3125
$ws.Range("K61").Value = 641.3125
$ws.Range("M61").Value = -429.3125
$ws.Range("H122").Value = 1133.5
$ws.Range("I122").Value = 1158.2106
$ws.Range("K122").Value = 3474.6318
$ws.Range("M122").Value = -1024.6318
$ws.Range("H132").Value = 2149.1904
$ws.Range("I132").Value = 1808.7333
$ws.Range("K132").Value = 5426.199900000001
$ws.Range("M132").Value = -2896.199900000001
$ws.Range("H136").Value = 641.3125
$ws.Range("I136").Value = 641.3125
$ws.Range("K136").Value = 1923.9375
$ws.Range("M136").Value = 626.0625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 45456060
$ws.Range("I99").Value = 50001468
$ws.Range("K99").Value = 50001468
$ws.Range("M99").Value = -49999970
$ws.Range("H107").Value = 1104.6842
$ws.Range("I107").Value = 752.5
$ws.Range("K107").Value = 752.5
$ws.Range("M107").Value = 1167.5
$ws.Range("H134").Value = 6470.9165
$ws.Range("I134").Value = 1016.8333
$ws.Range("J134").Value = 22833.166
$ws.Range("K134").Value = 3050.4999
$ws.Range("L134").Value = 68499.49800000001
$ws.Range("M134").Value = -515.4998999999998
$ws.Range("N134").Value = -73569.49800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 438.54544
$ws.Range("J22").Value = 950
$ws.Range("L22").Value = 950
$ws.Range("N22").Value = -1650
$ws.Range("H58").Value = 1449.8334
$ws.Range("I58").Value = 1264.5
$ws.Range("K58").Value = 1264.5
$ws.Range("M58").Value = -1061.5
$ws.Range("H132").Value = 8548.723
$ws.Range("I132").Value = 13462.889
$ws.Range("J132").Value = 3634.5557
$ws.Range("K132").Value = 40388.667
$ws.Range("L132").Value = 10903.6671
$ws.Range("M132").Value = -37858.667
$ws.Range("N132").Value = -15963.6671
$ws.Range("H134").Value = 1768
$ws.Range("I134").Value = 1793.8334
$ws.Range("K134").Value = 5381.5002
$ws.Range("M134").Value = -2846.5002
$ws.Range("H136").Value = 1449.8334
$ws.Range("I136").Value = 1264.5
$ws.Range("K136").Value = 3793.5
$ws.Range("M136").Value = -1243.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 16130276
$ws.Range("J131").Value = 1290.4576
$ws.Range("L131").Value = 3871.3728
$ws.Range("N131").Value = -13951.3728
$ws.Range("H134").Value = 4357.0435
$ws.Range("I134").Value = 2179.7778
$ws.Range("J134").Value = 5756.7144
$ws.Range("K134").Value = 6539.3334
$ws.Range("L134").Value = 17270.1432
$ws.Range("M134").Value = -1469.3334
$ws.Range("N134").Value = -27410.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 4973272.5
$ws.Range("I12").Value = 4770600
$ws.Range("J12").Value = 7000000
$ws.Range("K12").Value = 4770600
$ws.Range("L12").Value = 7000000
$ws.Range("M12").Value = -4770460
$ws.Range("N12").Value = -7000280
$ws.Range("H102").Value = 1833.8108
$ws.Range("I102").Value = 1845.7646
$ws.Range("K102").Value = 1845.7646
$ws.Range("M102").Value = -223.7646
$ws.Range("H122").Value = 1638.2667
$ws.Range("J122").Value = 1055.4
$ws.Range("L122").Value = 3166.2
$ws.Range("N122").Value = -8066.200000000001
$ws.Range("H126").Value = 2395
$ws.Range("I126").Value = 2348.3333
$ws.Range("K126").Value = 7044.999899999999
$ws.Range("M126").Value = -4574.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 551
$ws.Range("I16").Value = 537.8
$ws.Range("K16").Value = 537.8
$ws.Range("M16").Value = -367.8
$ws.Range("H40").Value = 2727.182
$ws.Range("I40").Value = 2428.8572
$ws.Range("K40").Value = 2428.8572
$ws.Range("M40").Value = -2292.8572
$ws.Range("H46").Value = 3000
$ws.Range("J46").Value = 3000
$ws.Range("L46").Value = 3000
$ws.Range("N46").Value = -3376
$ws.Range("H61").Value = 1551.8125
$ws.Range("I61").Value = 1540.4445
$ws.Range("J61").Value = 1566.4286
$ws.Range("K61").Value = 1540.4445
$ws.Range("L61").Value = 1566.4286
$ws.Range("M61").Value = -1338.4445
$ws.Range("N61").Value = -1970.4286
$ws.Range("H113").Value = 1551.8125
$ws.Range("I113").Value = 1540.4445
$ws.Range("J113").Value = 1566.4286
$ws.Range("K113").Value = 1540.4445
$ws.Range("L113").Value = 1566.4286
$ws.Range("M113").Value = 629.5554999999999
$ws.Range("N113").Value = -5906.4286
$ws.Range("H122").Value = 17865852
$ws.Range("I122").Value = 20843126
$ws.Range("J122").Value = 2205
$ws.Range("K122").Value = 62529378
$ws.Range("L122").Value = 6615
$ws.Range("M122").Value = -62526928
$ws.Range("N122").Value = -11515
$ws.Range("H136").Value = 5639.2856
$ws.Range("I136").Value = 7057.5
$ws.Range("K136").Value = 21172.5
$ws.Range("M136").Value = -18622.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9631975
$ws.Range("I122").Value = 12383479
$ws.Range("J122").Value = 1713.1666
$ws.Range("K122").Value = 37150437
$ws.Range("L122").Value = 5139.4998
$ws.Range("M122").Value = -37147987
$ws.Range("N122").Value = -10039.4998
$ws.Range("H132").Value = 4158.68
$ws.Range("I132").Value = 3606.2942
$ws.Range("K132").Value = 10818.8826
$ws.Range("M132").Value = -8288.882599999999
